$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet lists "missing items"; this update adds 4 new product
# rows (KETOLAC, VOLTAREN, "الويز كبير بالاجنحه" and "حلق") into the
# table, pushing the existing TORSERETIC and "كالونا" rows down, and
# pushing the totals row + footer row down as well. The grand total
# and the generation timestamp in the footer are refreshed too.
# ------------------------------------------------------------------

# --- Step 1: relocate the rows that already exist, working from the
#     bottom up so we never overwrite data before it has been copied.

# footer (row 13 -> row 17)
$ws.Range("A13:Q13").Copy($ws.Range("A17:Q17"))
# totals row (row 12 -> row 16)
$ws.Range("A12:Q12").Copy($ws.Range("A16:Q16"))
# "كالونا" row (row 11 -> row 15)
$ws.Range("A11:Q11").Copy($ws.Range("A15:Q15"))
# TORSERETIC row (row 10 -> row 11)
$ws.Range("A10:Q10").Copy($ws.Range("A11:Q11"))

# --- Step 2: build the 4 brand new rows (10, 12, 13, 14) by cloning
#     the formatting/merge layout of an existing data row, then
#     overwrite the cell values.

$ws.Range("A9:Q9").Copy($ws.Range("A10:Q10"))
$ws.Range("A9:Q9").Copy($ws.Range("A12:Q12"))
$ws.Range("A9:Q9").Copy($ws.Range("A13:Q13"))
$ws.Range("A9:Q9").Copy($ws.Range("A14:Q14"))

# --- Step 3: fill in the values for the new / shifted data rows.

# Row 10: KETOLAC 10MG 20 TAB
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "KETOLAC 10MG 20 TAB"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 1
$ws.Range("N10").Value = "38.00"
$ws.Range("P10").Value = "38.0000"
$ws.Range("Q10").Value = "1:0"
$ws.Rows("10:10").RowHeight = 24.75

# Row 11: TORSERETIC 100MG 30 TABS. (values unchanged, just renumbered)
$ws.Range("A11").Value = 5
$ws.Rows("11:11").RowHeight = 25.5

# Row 12: VOLTAREN 50MG 20 TAB.
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "VOLTAREN 50MG 20 TAB."
$ws.Range("H12").Value = "0:1"
$ws.Range("L12").Value = 1
$ws.Range("N12").Value = "48.00"
$ws.Range("P12").Value = "24.0000"
$ws.Range("Q12").Value = "0:1"
$ws.Rows("12:12").RowHeight = 25.5

# Row 13: الويز كبير بالاجنحه
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "الويز كبير بالاجنحه"
$ws.Range("H13").Value = "16:0"
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = "60.00"
$ws.Range("P13").Value = "60.0000"
$ws.Range("Q13").Value = "1:0"
$ws.Rows("13:13").RowHeight = 24.75

# Row 14: حلق
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "حلق"
$ws.Range("H14").Value = "22:0"
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = "10.00"
$ws.Range("P14").Value = "10.0000"
$ws.Range("Q14").Value = "1:0"
$ws.Rows("14:14").RowHeight = 25.5

# Row 15: كالونا (values unchanged, just renumbered)
$ws.Range("A15").Value = 9
$ws.Rows("15:15").RowHeight = 24.75

# --- Step 4: refresh the grand total and footer rows.

$ws.Range("P16").Value = 392.13
$ws.Rows("16:16").RowHeight = 25.5

$ws.Range("A17").Value = "Wednesday, 20 August, 2025 10:40 AM"
$ws.Rows("17:17").RowHeight = 16.5
